# daily auto push: 2026-01-31 02:46 UTC
# Insert a new data row right after the existing 2026/01/31 entries (row 743),
# pushing the rest of the log down by one row, and fill in the new entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 744 currently holds 2026/12/29 ... ; insert a new blank row above it so
# everything from 744 downward shifts to 745 downward.
$ws.Rows.Item(744).Insert()

# Force the date-looking text to stay literal text (not auto-converted to a
# date serial number), then drop the temporary "Text" number format so the
# cell ends up with the same default (no explicit style) as its neighbours.
$ws.Cells.Item(744, 1).NumberFormat = "@"
$ws.Cells.Item(744, 1).Value = "2026/01/31"
$ws.Cells.Item(744, 1).ClearFormats()

$ws.Cells.Item(744, 2).Value = "土"
$ws.Cells.Item(744, 3).Value = 8
$ws.Cells.Item(744, 4).Value = 201
